$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.522.18"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").Value = "1.877.51"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.36"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0704"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "2.147.58"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "12.11"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.10%  "
$ws.Range("D14").Value = "1.886.63"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.687"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").Value = "35.561.26"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "71.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").Value = "0.0₃0807"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("E20").Value = "  +0.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  +2.00%  "
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +26.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.29"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0565"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.01"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("E33").Value = "  +22.78%  "
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  +11.77%  "
$ws.Range("E36").Value = "  +5.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("E39").Value = "  +4.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "90.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").Value = "1.356.05"
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0591"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +11.20%  "
$ws.Range("E44").Value = "  +3.49%  "
$ws.Range("B45").Value = "Gas"
$ws.Range("C45").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +43.39%  "
$ws.Range("B46").Value = "HuobiToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +33.77%  "
$ws.Range("D50").Value = "2.066.12"
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0686"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.14%  "
